# Weekly update: insert a new week's worth of "Piña" price rows
# (4 rows, one per Calidad: Especial/Primera/Segunda/Tercera) at the top
# of the most-recent block, pushing all existing rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 990 (shifts old rows 990:1022 -> 994:1026)
$ws.Range("990:993").Insert()

# New data for the inserted rows. Non-varying columns (A,B,C,E,F,G,H,I,J,K,R)
# are identical to the rest of the sheet for this product/market.
$newRows = @(
    @{ Row = 990; L = "Especial"; N = 21500; O = 22000; P = 21750; Q = "$/caja 10 unidades"; S = 2175; T = 10 },
    @{ Row = 991; L = "Primera";  N = 21500; O = 22000; P = 21750; Q = "$/caja 12 unidades"; S = 1812; T = 12 },
    @{ Row = 992; L = "Segunda";  N = 21000; O = 22000; P = 21500; Q = "$/caja 14 unidades"; S = 1536; T = 14 },
    @{ Row = 993; L = "Tercera";  N = 21500; O = 22000; P = 21750; Q = "$/caja 16 unidades"; S = 1359; T = 16 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 8
    $ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44747
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100108
    $ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($row, 9).Value = 100108005
    $ws.Cells.Item($row, 10).Value = "Piña"
    $ws.Cells.Item($row, 11).Value = "Caramelo"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = 216
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Ecuador"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
